# Generate Report for Handback
# Applies the "handback" localization-status update to the zh-cn and de-de
# report sheets: marks the status as synced, fills in the Latest Target
# File / Latest Handback File / Latest Handback DateTime columns, and
# widens a couple of columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5b731a78c4e382ed5b2db7152c30f59c1bb99861/e2e"
$mdName = "e970ca74-c046-45e8-af28-eb53490dd6c6.md"
$mdUrl = "$githubBase/$mdName"

$statusText = "Handed back: in sync with en-US"

# ColumnWidth inputs chosen so the stored/exported column width lands on
# the same pixel-quantized width for every column that needs to widen.
$wideColWidth = 29.16666667   # -> stored width ~30 ("Status" column)
$maxColWidth  = 39.16666667   # -> stored width 40 (Latest Target File / Latest Handback File)

$zhXlf = "e970ca74-c046-45e8-af28-eb53490dd6c6.a8e3dabb0a8d509bf793b86461d1878a7ad44a81.zh-cn.xlf"
$deXlf = "e970ca74-c046-45e8-af28-eb53490dd6c6.a8e3dabb0a8d509bf793b86461d1878a7ad44a81.de-de.xlf"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZh.Columns.Item(9).ColumnWidth = $maxColWidth
$wsZh.Columns.Item(10).ColumnWidth = $maxColWidth

$wsZh.Cells.Item(2, 3).Value = $statusText
$wsZh.Cells.Item(3, 3).Value = $statusText

$wsZh.Cells.Item(2, 9).Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(2, 9), $mdUrl, "", "", $mdName)
$wsZh.Cells.Item(3, 9).Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 9), $mdUrl, "", "", $mdName)

$wsZh.Cells.Item(2, 10).Value = $zhXlf
$wsZh.Cells.Item(3, 10).Value = $zhXlf

$wsZh.Cells.Item(2, 11).Value = "2016-09-05 15:17:55"
$wsZh.Cells.Item(3, 11).Value = "2016-09-05 15:17:55"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDe.Columns.Item(9).ColumnWidth = $maxColWidth
$wsDe.Columns.Item(10).ColumnWidth = $maxColWidth

$wsDe.Cells.Item(2, 3).Value = $statusText
$wsDe.Cells.Item(3, 3).Value = $statusText

$wsDe.Cells.Item(2, 9).Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(2, 9), $mdUrl, "", "", $mdName)
$wsDe.Cells.Item(3, 9).Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 9), $mdUrl, "", "", $mdName)

$wsDe.Cells.Item(2, 10).Value = $deXlf
$wsDe.Cells.Item(3, 10).Value = $deXlf

$wsDe.Cells.Item(2, 11).Value = "2016-09-05 15:18:08"
$wsDe.Cells.Item(3, 11).Value = "2016-09-05 15:18:08"

# ---------------------------------------------------------------------
# Overview sheet — its zh-cn / de-de columns (E, F) mirror the same
# "Status" shared string as the per-language sheets, so they pick up the
# new status text as well; widen those columns to fit it.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

$wsOverview.Cells.Item(2, 5).Value = $statusText
$wsOverview.Cells.Item(2, 6).Value = $statusText
$wsOverview.Cells.Item(3, 5).Value = $statusText
$wsOverview.Cells.Item(3, 6).Value = $statusText
